$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows 87-105: file name (col A), title (col B), JSON-builder formula (col C)
$rows = @(
    @{ R = 87;  A = "airplanejivetalkersfoodsickwithjunecleaver.mp4"; B = "Airplane Movie -Jive Talkers Sick with June Cleaver"; Pct = $true },
    @{ R = 88;  A = "airplanejivetalkers.mp4"; B = "Airplane Movie -Jive Talkers"; Pct = $true },
    @{ R = 89;  A = "airplanejivetalkersorderfood.mp4"; B = "Airplane Movie -Jive Talkers Order Food"; Pct = $true },
    @{ R = 90;  A = "01 Here To Fall.mp3"; B = "Here to Fall"; Pct = $true },
    @{ R = 91;  A = "02 Avalon Or Someone Very Similar.mp3"; B = "Avalon Or Someone Very Similar"; Pct = $true },
    @{ R = 92;  A = "03 By Twos.mp3"; B = "By Twos"; Pct = $true },
    @{ R = 93;  A = "04 Nothing To Hide.mp3"; B = "Nothing to Hide"; Pct = $true },
    @{ R = 94;  A = "05 Periodically Triple Or Double.mp3"; B = "Periodically Triple Or Double"; Pct = $true },
    @{ R = 95;  A = "06 If Its True.mp3"; B = "If Its True"; Pct = $true },
    @{ R = 96;  A = "07 Im On My Way.mp3"; B = "Im On My Way"; Pct = $true },
    @{ R = 97;  A = "08 When Its Dark.mp3"; B = "When Its Dark"; Pct = $true },
    @{ R = 98;  A = "09 All Your Secrets.mp3"; B = "All Your Secrets"; Pct = $true },
    @{ R = 99;  A = "10 More Stars Than There Are In Heaven.mp3"; B = "More Stars Than There Are In Heaven"; Pct = $true },
    @{ R = 100; A = "11 The Fireside.mp3"; B = "The Fireside"; Pct = $true },
    @{ R = 101; A = "IMG_4599.JPG"; B = "Kitteh Mehowh"; Pct = $true },
    @{ R = 102; A = "IMG_4659.JPG"; B = "Daffodil"; Pct = $true },
    @{ R = 103; A = "IMG_4675.JPG"; B = "Good Luck Kitteh!"; Pct = $false },
    @{ R = 104; A = "IMG_4575.JPG"; B = "Good Boy!"; Pct = $false },
    @{ R = 105; A = "IMG_4691.JPG"; B = "Leaf Blower"; Pct = $false }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    if ($row.Pct) {
        $ws.Cells.Item($r, 2).NumberFormat = "0%"
    }
    $ws.Cells.Item($r, 3).Formula = "=""{'file':'""&A$r&""','title':'""&B$r&""'},"""
}

# Column A width
$ws.Columns.Item(1).ColumnWidth = 40.666666666666664

# Scroll / selection to match final view
$ws.Range("C101:C105").Select() | Out-Null
